$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update shared strings used by the "临时功能" column for FSK control rows
$ws.Range("D3").Value = "PWM_CH_1"
$ws.Range("D4").Value = "FSK控制输出,PWM_CH_4"

# Update the saved cell selection on the sheet view
$ws.Range("D5").Select()
